# Insert two new rows at row 147 (shifting settingsTabs.Time etc. down by 2)
# and populate them with the new crossPlot.sidebar.* translation keys.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()
$ws.Rows.Item(147).Insert()

$ws.Range("A147").Value = "crossPlot.sidebar.no-markers"
$ws.Range("B147").Value = " No Markers Available (Please Add a Marker by Clicking the Chart in the Main View)"

$ws.Range("A148").Value = "crossPlot.sidebar.no-models"
$ws.Range("B148").Value = "No Models Available (Please Add a Model by Clicking the Chart in the Main View)"
